# feat: add 2022-Q4 data
#
# 1) Insert a brand-new "2022-Q4" worksheet right after "总计" (before "2022-Q3"),
#    populated with the fund-holding detail rows for that quarter.
# 2) Insert a new summary row for 2022-Q4 at the top of the "总计" sheet,
#    shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q4" sheet, positioned before "2022-Q3"
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row (bold, bordered, centered - same look as the other quarter sheets)
$headerRange = $q4.Range("B1:H1")
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows - D/E/F/G are text in the source data, so force "@" text format
# before assignment so Excel doesn't coerce them to numbers.
$q4.Range("A2").Value = 0
$q4.Range("A2").Font.Bold = $true
$q4.Range("A2").Borders.LineStyle = 1
$q4.Range("A2").HorizontalAlignment = -4108
$q4.Range("A2").VerticalAlignment = -4160
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "000926"
$q4.Range("C2").Value = "中信建投睿信灵活配置混合A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.10"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "83.25"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "1.05"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0010"
$q4.Range("H2").Value = 5

$q4.Range("A3").Value = 1
$q4.Range("A3").Font.Bold = $true
$q4.Range("A3").Borders.LineStyle = 1
$q4.Range("A3").HorizontalAlignment = -4108
$q4.Range("A3").VerticalAlignment = -4160
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "004676"
$q4.Range("C3").Value = "中信建投睿信灵活配置混合C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.03"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "83.25"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "1.05"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0003"
$q4.Range("H3").Value = 5

# ---------------------------------------------------------------------------
# 2. Insert the 2022-Q4 summary row into "总计", above the existing 2022-Q3 row
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("A2").Font.Bold = $true
$total.Range("A2").Borders.LineStyle = 1
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

Write-Output "2022-Q4 sheet added and summary sheet updated"
